$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "91.596.68"
$ws.Range("E2").Value = "  +5.30%  "

# Row 3
$ws.Range("D3").Value = "3.300.31"
$ws.Range("E3").Value = "  +1.03%  "

# Row 4
$ws.Range("E4").Value = "  +0.22%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.48"
$ws.Range("E5").Value = "  +1.27%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "627.66"
$ws.Range("E6").Value = "  -0.22%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.418"
$ws.Range("E7").Value = "  +11.15%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.709"
$ws.Range("E8").Value = "  +2.23%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  +0.03%  "

# Row 10
$ws.Range("D10").Value = "3.295.98"
$ws.Range("E10").Value = "  +1.14%  "

# Row 11
$ws.Range("E11").Value = "  +2.35%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000266"
$ws.Range("E12").Value = "  +3.03%  "

# Row 13
$ws.Range("E13").Value = "  +1.08%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.41"
$ws.Range("E14").Value = "  +0.50%  "

# Row 15
$ws.Range("D15").Value = "3.909.85"
$ws.Range("E15").Value = "  +1.10%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.39"
$ws.Range("E16").Value = "  +1.49%  "

# Row 17
$ws.Range("D17").Value = "91.185.29"
$ws.Range("E17").Value = "  +4.95%  "

# Row 18
$ws.Range("D18").Value = "3.293.28"
$ws.Range("E18").Value = "  +0.52%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.24"
$ws.Range("E19").Value = "  +5.47%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.15"
$ws.Range("E20").Value = "  +0.58%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "435.66"
$ws.Range("E21").Value = "  +0.59%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.04"
$ws.Range("E22").Value = "  +1.29%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.34"
$ws.Range("E23").Value = "  +0.11%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000185"
$ws.Range("E24").Value = "  +42.83%  "

# Row 25
$ws.Range("E25").Value = "  +5.94%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.21"
$ws.Range("E26").Value = "  -2.31%  "

# Row 27
$ws.Range("D27").Value = "3.471.37"
$ws.Range("E27").Value = "  +0.88%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "76.43"
$ws.Range("E28").Value = "  +0.42%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.14%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.180"
$ws.Range("E30").Value = "  +3.68%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.40%  "

# Row 32
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "561.16"
$ws.Range("E32").Value = "  +2.86%  "

# Row 33
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.70"
$ws.Range("E33").Value = "  -1.30%  "

# Row 34
$ws.Range("E34").Value = "  +6.46%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.72"
$ws.Range("E35").Value = "  +26.85%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.34"
$ws.Range("E36").Value = "  -6.20%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.92"
$ws.Range("E37").Value = "  -0.96%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.62"
$ws.Range("E38").Value = "  +0.45%  "

# Row 39
$ws.Range("E39").Value = "  -2.89%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.37"
$ws.Range("E40").Value = "  +3.61%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.47%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.397"
$ws.Range("E42").Value = "  +0.31%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.99"
$ws.Range("E43").Value = "  -0.42%  "

# Row 44
$ws.Range("E44").Value = "  +0.05%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "183.48"
$ws.Range("E45").Value = "  +2.22%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "149.00"
$ws.Range("E46").Value = "  -5.28%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.130"
$ws.Range("E47").Value = "  +5.80%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "43.94"
$ws.Range("E48").Value = "  -1.07%  "

# Row 49
$ws.Range("E49").Value = "  -1.10%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.91"
$ws.Range("E50").Value = "  +7.21%  "

# Row 51
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.629"
$ws.Range("E51").Value = "  +0.40%  "
